$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.373.24"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.885.41"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "'245.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.01%  "
$ws.Range("D6").Value = "'0.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "'43.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.12%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("D10").Value = "'53.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "'0.0972"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'13.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "'0.765"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("E16").Value = "  -1.65%  "
$ws.Range("D17").Value = "1.885.81"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "35.525.11"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "'73.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "'244.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'12.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("E23").Value = "  -1.48%  "
$ws.Range("E24").Value = "  +9.06%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("D27").Value = "'165.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "'18.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "4.128.44"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +9.74%  "
$ws.Range("E33").Value = "  -1.14%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").Value = "'1.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.67%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").Value = "'0.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").Value = "'0.0701"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.33%  "
$ws.Range("D41").Value = "'17.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'0.0218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'96.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.87%  "
$ws.Range("E44").Value = "  -2.11%  "
$ws.Range("D45").Value = "1.307.66"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "'0.0797"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.91%  "
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "'12.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -5.11%  "
